# Update cryptos list values (price + 1h volume change) per Oct 6 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.983.40'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.418.98'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''563.37'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '''142.73'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.531'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '''5.20'
$ws.Range("E11").Value = '  -3.94%  '
$ws.Range("D12").Value = '''0.350'
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Value = '''25.91'
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '2.855.75'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '61.881.09'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '2.406.37'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '''11.32'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = '''324.08'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("D21").Value = '''6.82'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''66.69'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").Value = '''8.70'
$ws.Range("E25").Value = '  -3.57%  '
$ws.Range("D26").Value = '''549.47'
$ws.Range("E26").Value = '  -7.58%  '
$ws.Range("D27").Value = '2.538.29'
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '0.0₃0931'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '''8.16'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  -4.23%  '
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("E34").Value = '  -4.03%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").Value = '''153.30'
$ws.Range("E38").Value = '  +1.26%  '
$ws.Range("D39").Value = '''5.41'
$ws.Range("E39").Value = '  -5.38%  '
$ws.Range("D40").Value = '''18.55'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").Value = '''1.80'
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").Value = '''146.65'
$ws.Range("E43").Value = '  -2.99%  '
$ws.Range("E44").Value = '  -6.21%  '
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = '''0.0527'
$ws.Range("E46").Value = '  -2.75%  '
$ws.Range("D47").Value = '''19.80'
$ws.Range("E47").Value = '  -2.47%  '
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("E51").Value = '  +0.62%  '
